$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# D-column values are forced to remain text (matching the original inline-string
# cell type) even when they look like plain numbers (e.g. "235.11"), by temporarily
# applying a text number format and then restoring the default "Normal" style so
# no extraneous formatting change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.427.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.069.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.28'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.54%  '

$ws.Range("E9").Value = '  +2.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0773'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.61%  '

$ws.Range("E11").Value = '  +0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.373.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("E13").Value = '  -1.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.777'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.069.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.330.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0818'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("E25").Value = '  -1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.70%  '

$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("E28").Value = '  -6.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.117'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("E33").Value = '  -0.86%  '

$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("E35").Value = '  -2.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.91%  '

$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("E39").Value = '  -4.17%  '

$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  -3.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.483.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.64%  '

$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("E45").Value = '  -0.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.42%  '

$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.59%  '

$ws.Range("E49").Value = '  -1.07%  '

$ws.Range("E50").Value = '  +0.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.14%  '
